$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Propagate existing cell formatting to the new cells so the
#        header row (style 1) and body rows (style 2) match the
#        original look-and-feel before we fill in values. ---

# New header cells F1, G1 should look like the existing header cell E1.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("E1").Copy($ws.Range("G1"))

# New body cells F2, G2 should look like the existing body cell E2.
$ws.Range("E2").Copy($ws.Range("F2"))
$ws.Range("E2").Copy($ws.Range("G2"))

# Brand new rows 3 and 4 should look like row 2 across all 7 columns.
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("C2").Copy($ws.Range("C3"))
$ws.Range("D2").Copy($ws.Range("D3"))
$ws.Range("E2").Copy($ws.Range("E3"))
$ws.Range("F2").Copy($ws.Range("F3"))
$ws.Range("G2").Copy($ws.Range("G3"))

$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("B2").Copy($ws.Range("B4"))
$ws.Range("C2").Copy($ws.Range("C4"))
$ws.Range("D2").Copy($ws.Range("D4"))
$ws.Range("E2").Copy($ws.Range("E4"))
$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("G2").Copy($ws.Range("G4"))

# --- 2. Fill in the real values for the Contacts test-case grid. ---

# Row 1 - headers
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Middle Name"
$ws.Range("D1").Value = "Last Name"
$ws.Range("E1").Value = "Suffix"
$ws.Range("F1").Value = "File Path"
$ws.Range("G1").Value = "RunMode"

# Row 2
$ws.Range("A2").Value = "Mr."
$ws.Range("B2").Value = "Ashish"
$ws.Range("C2").Value = "Rajesh"
$ws.Range("D2").Value = "Goyal"
$ws.Range("E2").Value = "Esq."
$ws.Range("F2").Value = "D:\ashishProject\freecrmproject\TestData\image1.jfif"
$ws.Range("G2").Value = "Y"

# Row 3
$ws.Range("A3").Value = "Mrs."
$ws.Range("B3").Value = "Anisha"
$ws.Range("C3").Value = "Rajesh"
$ws.Range("D3").Value = "Goyal"
$ws.Range("E3").Value = "II"
$ws.Range("F3").Value = "D:\ashishProject\freecrmproject\TestData\image2.jfif"
$ws.Range("G3").Value = "Y"

# Row 4
$ws.Range("A4").Value = "Mrs."
$ws.Range("B4").Value = "Snehal"
$ws.Range("C4").Value = "Ashish"
$ws.Range("D4").Value = "Goyal"
$ws.Range("E4").Value = "III"
$ws.Range("F4").Value = "D:\ashishProject\freecrmproject\TestData\image1.jfif"
$ws.Range("G4").Value = "Y"

# --- 3. Column widths for the two new columns (E already existed with
#        the bestFit width carried over; widen E slightly and size F to
#        comfortably fit the long file-path strings). ---
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 50.8

# --- 4. Rows 3 & 4 wrap across two lines (Suffix/File Path columns),
#        so give them the taller row height. ---
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8

# --- 5. Selection, matching the saved cursor position. ---
$null = $ws.Range("F12").Select()
